{"js": "// Remove the \"question\" block paragraphs (from the \"Arist\u00f3teles, S\u00f3crates e\n// Plat\u00e3o...\" paragraph through the blank paragraph right before \"Resolu\u00e7\u00e3o\")\n// and relocate the \"_GoBack\" bookmark from its old spot (right after the\n// \"Matematica...\" text run, near the end of the document) to the very start\n// of the \"Resolu\u00e7\u00e3o\" paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst texts = paragraphs.items.map((p) => p.text);\n\n// Locate the start of the block to remove: the paragraph that begins the\n// \"Arist\u00f3teles, S\u00f3crates e Plat\u00e3o...\" question text.\nconst blockStartText =\n  \"Arist\u00f3teles, S\u00f3crates e Plat\u00e3o estavam disputando uma competi\u00e7\u00e3o de perguntas nas disciplinas de F\u00edsica, Matem\u00e1tica e Qu\u00edmica. Cada um obteve um primeiro lugar, um segundo lugar e um terceiro lugar.\";\n// Locate the end of the block to remove: the paragraph right after \"Informe\n// quem ficou...\" (a blank paragraph), i.e. the paragraph immediately before\n// \"Resolu\u00e7\u00e3o\".\nconst resolucaoText = \"Resolu\u00e7\u00e3o\";\n\nlet startIdx = -1;\nlet resolucaoIdx = -1;\nfor (let i = 0; i < texts.length; i++) {\n  if (startIdx === -1 && texts[i] === blockStartText) {\n    startIdx = i;\n  }\n  if (texts[i] === resolucaoText) {\n    resolucaoIdx = i;\n    break;\n  }\n}\n\nif (startIdx === -1 || resolucaoIdx === -1 || resolucaoIdx <= startIdx) {\n  throw new Error(\n    \"Could not locate the expected paragraphs (startIdx=\" +\n      startIdx +\n      \", resolucaoIdx=\" +\n      resolucaoIdx +\n      \").\"\n  );\n}\n\n// Delete every paragraph from the question block start up to (but not\n// including) the \"Resolu\u00e7\u00e3o\" paragraph.\nfor (let i = resolucaoIdx - 1; i >= startIdx; i--) {\n  paragraphs.items[i].delete();\n}\nawait context.sync();\n\n// Move the \"_GoBack\" bookmark to the start of the (now immediately\n// following) \"Resolu\u00e7\u00e3o\" paragraph.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst refreshedParagraphs = body.paragraphs;\nrefreshedParagraphs.load(\"text\");\nawait context.sync();\n\nlet newResolucaoIdx = -1;\nfor (let i = 0; i < refreshedParagraphs.items.length; i++) {\n  if (refreshedParagraphs.items[i].text === resolucaoText) {\n    newResolucaoIdx = i;\n    break;\n  }\n}\nif (newResolucaoIdx === -1) {\n  throw new Error(\"Could not relocate the 'Resolu\u00e7\u00e3o' paragraph after deletion.\");\n}\n\nconst resolucaoStart = refreshedParagraphs.items[newResolucaoIdx].getRange(\"Start\");\nresolucaoStart.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Remove the \"question\" block paragraphs (from the \"Arist\u00f3teles, S\u00f3crates e\n# Plat\u00e3o...\" paragraph through the blank paragraph right before \"Resolu\u00e7\u00e3o\")\n# and relocate the \"_GoBack\" bookmark from its old spot (right after the\n# \"Matematica...\" text run, near the end of the document) to the very start\n# of the \"Resolu\u00e7\u00e3o\" paragraph.\n\n$d = $word.ActiveDocument\n\n$startMarker = \"Arist\u00f3teles, S\u00f3crates e Plat\u00e3o estavam disputando\"\n$resolucaoMarker = \"Resolu\u00e7\u00e3o\"\n\n$startIndex = -1\n$resolucaoIndex = -1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs($i).Range.Text\n    if ($startIndex -eq -1 -and $t -like \"$startMarker*\") {\n        $startIndex = $i\n    }\n    if ($t -like \"$resolucaoMarker*\") {\n        $resolucaoIndex = $i\n        break\n    }\n}\n\nif ($startIndex -eq -1 -or $resolucaoIndex -eq -1 -or $resolucaoIndex -le $startIndex) {\n    throw \"Could not locate the expected paragraphs (startIndex=$startIndex, resolucaoIndex=$resolucaoIndex).\"\n}\n\n$endIndex = $resolucaoIndex - 1\n\n$startRange = $d.Paragraphs($startIndex).Range\n$endRange = $d.Paragraphs($endIndex).Range\n$blockRange = $d.Range($startRange.Start, $endRange.End)\n$blockRange.Delete()\n\n# Move the \"_GoBack\" bookmark to the start of the (now immediately\n# following) \"Resolu\u00e7\u00e3o\" paragraph.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs($i).Range.Text\n    if ($t -like \"$resolucaoMarker*\") {\n        $target = $d.Paragraphs($i).Range.Duplicate\n        $target.Collapse(1)  # wdCollapseStart\n        $d.Bookmarks.Add(\"_GoBack\", $target)\n        break\n    }\n}\n"}
